$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "outside"
$ws.Range("B3").Value = "home"
